$d = $word.ActiveDocument

# --- 1. Remove the existing "_GoBack" bookmark -----------------------------
# It currently sits right after "Le chiffrement ElGamal (1984) ", before the
# "est une variante..." run. It will be re-created at the end of the new
# paragraph added below.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- 2. Find the paragraph "-d'une clé publique égale à g^s mod p." --------
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*clé publique égale à g^s mod p.*") {
        $targetPara = $p
        break
    }
}
if ($null -eq $targetPara) {
    throw "Could not find the '-d'une clé publique égale à g^s mod p.' paragraph"
}

# --- 3. Insert a brand-new paragraph right after it and fill it in ---------
$targetPara.Range.InsertParagraphAfter()
$newPara = $targetPara.Next()
$newPara.Range.Text = "La sécurité repose sur le fait qu’il est « difficile » de calculer s à partir de g^s dans Z/pZ. Ainsi la connaissance de la clé publique de B ne permet pas d’obtenir sa clé privée."

# --- 4. Re-create "_GoBack" collapsed at the end of the new paragraph's ----
#        text (i.e. right after the sentence, before the paragraph mark).
#
# Note: adding a bookmark directly with a collapsed range sitting exactly at
# "paragraph.End - 1" trips a boundary bug in this host's Bookmarks.Add, so
# a trailing marker character is appended first to move the desired offset
# away from that exact boundary, the bookmark is added there, and then the
# marker is deleted again - the now-collapsed bookmark stays in place.
$newPara.Range.InsertAfter("#")
$markerPos = $newPara.Range.End - 2
$bmRange = $d.Range($markerPos, $markerPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
$d.Range($markerPos, $markerPos + 1).Delete()

Write-Output "done"
